$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the title text (shared string used in A1)
$ws.Range("A1").Value = "Total Load Profile for All Users on 2023-01-20"

# Update the individual customer demand values in column B (rows 3-24)
$ws.Range("B3").Value = 0.6567142857
$ws.Range("B4").Value = 0.6577142857
$ws.Range("B5").Value = 0.6647142857
$ws.Range("B6").Value = 0.5467142857
$ws.Range("B7").Value = 0.2647142857
$ws.Range("B8").Value = 0.2887142857
$ws.Range("B9").Value = 0.3106345858
$ws.Range("B10").Value = 0.3197117632
$ws.Range("B11").Value = 0.4647965081
$ws.Range("B12").Value = 0.551
$ws.Range("B13").Value = 0.533
$ws.Range("B14").Value = 0.499
$ws.Range("B15").Value = 0.4965
$ws.Range("B16").Value = 0.5115000000000001
$ws.Range("B17").Value = 0.506
$ws.Range("B18").Value = 0.488
$ws.Range("B19").Value = 0.7458028329
$ws.Range("B20").Value = 1.4707985559
$ws.Range("B21").Value = 1.4083986112
$ws.Range("B22").Value = 1.414
$ws.Range("B23").Value = 0
$ws.Range("B24").Value = 0
